$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the breed filter from "Miniature Schnauzer" to "Miniature Pinscher" in the
# query text cells: B2 (Cases query), B3 (Samples query), B4 (Files query), and the
# shared stat query in column C (C2/C3/C4 all hold the same query text).
$cells = @("B2", "B3", "B4", "C2", "C3", "C4")
foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $text = $range.Value2
    $range.Value2 = $text -replace "Miniature Schnauzer", "Miniature Pinscher"
}

# Update the view state: scroll so row 4 is the top-left visible row, zoom to 70%,
# and move the active selection to B4.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 70
$null = $ws.Range("B4").Select()
